$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","I","J","L","M","N")

$data = @{
    2 = @(33.23047904178452, 19.86901821678549, 11.1547285587345, 9.269531559306516, 3.861643185877579, 46.17174936111569, 6.344284098028226, 15.00686773499833, 26.16377165705835, 25.63916579089601)
    3 = @(33.03940577600938, 19.59037697689681, 11.1751999113411, 9.258108059337692, 3.866860942541714, 46.13933565104431, 6.350943247711864, 15.02834488693622, 26.14950876479657, 25.67358032268842)
    4 = @(32.93094029205604, 19.42384398051451, 11.18933253883156, 9.250967766169358, 3.870227441466213, 46.12772438176076, 6.356207020055137, 15.04394504891497, 26.14640633426343, 25.69661449019479)
    5 = @(32.88899833142332, 19.35720664796712, 11.19548463183433, 9.248025951396945, 3.871640426794786, 46.12507461749587, 6.358646989444681, 15.05090833259212, 26.14656409532334, 25.70647898579796)
    6 = @(32.88217117949856, 19.34621794396421, 11.19652991016342, 9.247535517110933, 3.871877540087739, 46.12476026342038, 6.359069941252484, 15.05210116479036, 26.14667616287026, 25.70814581857626)
    7 = @(32.93036546339377, 19.42294021614732, 11.18941391731365, 9.250928222210288, 3.870246330790003, 46.12768022103821, 6.356238732879674, 15.0440365051234, 26.14640270468655, 25.69674559227938)
    8 = @(33.16277967720185, 19.77204704780868, 11.16146257905321, 9.265618720945959, 3.86340859213677, 46.1588509036461, 6.346335995217624, 15.01377181493319, 26.15768035060238, 25.65063638478398)
    9 = @(33.68721464398347, 20.48909676107139, 11.11905877624877, 9.293441850839827, 3.8512832083844, 46.28585284660496, 6.376104971456146, 14.97359963678973, 26.22463479307174, 25.57535500156725)
    10 = @(34.11203015828291, 21.03066829019326, 11.09547758653459, 9.31330730286672, 3.843145673405085, 46.41939340016732, 6.455893136271731, 14.95581437442559, 26.30107499111924, 25.52932424971554)
    11 = @(34.31330516125492, 21.27917230897687, 11.08639568912072, 9.322224782241733, 3.839608660773862, 46.48887340043949, 6.491330280368625, 14.95027690820797, 26.34173079222488, 25.5104094480333)
    12 = @(34.39062512526432, 21.37349333619802, 11.08319326866747, 9.325584825514094, 3.838292792789315, 46.51643639261152, 6.504624211988011, 14.94854749440462, 26.35796744403936, 25.50353902573977)
    13 = @(34.37392475573353, 21.35317133317862, 11.08387243866692, 9.324861922092524, 3.838575145145999, 46.51044457017399, 6.501766729605427, 14.948903604969, 26.35443326350881, 25.50500568158585)
    14 = @(34.31964455359009, 21.28692829682515, 11.08612747979715, 9.322501549909109, 3.839499933111334, 46.49111594635097, 6.492426502628039, 14.95012726121345, 26.3430497630725, 25.50983835209271)
    15 = @(34.28653821646365, 21.24637827163026, 11.08753958519062, 9.321053574984413, 3.840069450496759, 46.47943959427675, 6.486688966303173, 14.95092465509734, 26.3361864314129, 25.51283658675218)
    16 = @(34.09903340496713, 21.01446332042251, 11.0961042271182, 9.312722152342383, 3.843380126585545, 46.41502817743669, 6.453559794471239, 14.9562276701115, 26.29853602861586, 25.5306012091729)
    17 = @(33.98602393451534, 20.87267478746411, 11.1017798417412, 9.307580764394812, 3.845453201947317, 46.37774836548786, 6.433014818135735, 14.96013511005824, 26.27694233630682, 25.5420185657111)
    18 = @(33.92178245092359, 20.79132853059181, 11.10519916439833, 9.304612312605864, 3.846661101615747, 46.35712814093908, 6.421117155825313, 14.96262286409421, 26.26507626957693, 25.54877604973406)
    19 = @(33.90016323367547, 20.76382436090337, 11.1063834819082, 9.30360530077772, 3.847072746861323, 46.35028776801951, 6.417075002070121, 14.96350643345507, 26.26115391321004, 25.55109670877863)
    20 = @(33.99797584846117, 20.88774767651823, 11.10115963469954, 9.308129236837917, 3.845230914383658, 46.38163179647158, 6.435210235169297, 14.95969428352276, 26.27918371095618, 25.54078344282075)
    21 = @(34.33555849895207, 21.30638027435669, 11.0854586949642, 9.323195301856531, 3.839227663484591, 46.49675927240375, 6.495173370008768, 14.94975786731932, 26.34637058576484, 25.50841094096939)
    22 = @(34.56258283605634, 21.58121097326995, 11.07657690437745, 9.332944017506291, 3.835441214581444, 46.57929952531179, 6.533631309650087, 14.94540604175304, 26.39518134029433, 25.48895731331818)
    23 = @(34.4408484248409, 21.43444529009201, 11.08119101450893, 9.327749776696347, 3.837449633884402, 46.53457984199147, 6.513173106919097, 14.94753258648516, 26.36868353966168, 25.49918384719034)
    24 = @(33.99257011302929, 20.88093268746793, 11.10143954335897, 9.307881311717592, 3.84533136055111, 46.37987356653944, 6.434217955113747, 14.95989282958664, 26.27816867687573, 25.54134123926842)
    25 = @(33.53821932182072, 20.29215810311216, 11.12920075335383, 9.286017897038889, 3.854427229163647, 46.24443095865834, 6.345972869059222, 14.98240992943324, 26.20172702684741, 25.57535500156725)
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $cell = "$($cols[$i])$row"
        $ws.Range($cell).Value = $vals[$i]
    }
}
